# Update Thresholds and Results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 0.80250335769744807
$ws.Range("G2").Value = 0.91648351648351656
$ws.Range("H2").Value = 0.8717948717948717

# Row 6
$ws.Range("E6").Value = 0.72341598580802502
$ws.Range("G6").Value = 0.87816225646990409
$ws.Range("H6").Value = 0.75000000000000011

# Row 8
$ws.Range("E8").Value = 0.75012146732608598
$ws.Range("H8").Value = 0.72727272727272729

# Row 9
$ws.Range("E9").Value = 0.85676283368939932
$ws.Range("F9").Value = 0.98
$ws.Range("G9").Value = 0.9893617021276595
$ws.Range("H9").Value = 0.8571428571428571

# Row 10
$ws.Range("D10").Value = 0.998
$ws.Range("E10").Value = 0.7673065735892961
$ws.Range("F10").Value = 0.98
$ws.Range("G10").Value = 0.88365328679464805
$ws.Range("H10").Value = 0.77777777777777779

# Row 12
$ws.Range("D12").Value = 0.98000000000000009
$ws.Range("E12").Value = 0.7376717048657514
$ws.Range("F12").Value = 0.98
$ws.Range("H12").Value = 0.7142857142857143

# Row 13
$ws.Range("E13").Value = 0.92343088262881712
$ws.Range("F13").Value = 0.995
$ws.Range("G13").Value = 0.99742268041237114
$ws.Range("H13").Value = 0.92307692307692302

# Reset the sheet selection to A1 (author re-saved with no extra selection)
$ws.Range("A1").Select()
